# Excel reading utility: adds an "expectedUser" verification column plus
# webpage name/URL columns to the test-data sheet, then centers the data
# and resizes the columns to fit their new contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers (write webpageURL before webpageName so the shared-string
# table allocates the indices in the order the source workbook has them)
$ws.Range("C1").Value = "expectedUser"
$ws.Range("F1").Value = "webpageURL"
$ws.Range("E1").Value = "webpageName"

# New "expectedUser" column mirrors the username column for each row
$ws.Range("C2").Value = "admin"
$ws.Range("C3").Value = "Test001"
$ws.Range("C4").Value = "Test002"
$ws.Range("C5").Value = "Test003"
$ws.Range("C6").Value = "Test004"
$ws.Range("C7").Value = "Test005"

# Center-align the populated cells (this is how the authored sheet stores
# style index 1 -- applied to every cell that actually holds data)
$ws.Range("A1:C7").HorizontalAlignment = -4108
$ws.Range("E1:F1").HorizontalAlignment = -4108

# Resize columns to fit their (new) contents
$ws.Columns.Item(1).ColumnWidth = 8.333333333333334
$ws.Columns.Item(2).ColumnWidth = 8
$ws.Columns.Item(3).ColumnWidth = 11.333333333333334
$ws.Columns.Item(4).ColumnWidth = 7.833333333333333
$ws.Columns.Item(5).ColumnWidth = 12.5
$ws.Columns.Item(6).ColumnWidth = 10.666666666666666

# Leave the final selection where the author left it
$ws.Range("F7").Select() | Out-Null
